# Final re-run including fig: unmet social support
# Adds an "Associated variable" column to the "Changes to be made" sheet,
# recording which associated variable (DSD1_A1 / DSD1_A9) goes with each
# re-classified record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes to be made")

# Header for the new column D
$ws.Range("D2").Value = "Associated variable"
$ws.Range("D2").Font.Bold = $true

# Row 4 (ID 27, Schizophrénia) uses DSD1_A1; every other data row uses DSD1_A9
$ws.Range("D4").Value = "DSD1_A1"
$ws.Range("D3").Value = "DSD1_A9"
$ws.Range("D5").Value = "DSD1_A9"
$ws.Range("D6").Value = "DSD1_A9"
$ws.Range("D7").Value = "DSD1_A9"
$ws.Range("D8").Value = "DSD1_A9"
$ws.Range("D9").Value = "DSD1_A9"
$ws.Range("D10").Value = "DSD1_A9"

$ws.Activate()
